# Refactored the test to remove individual products, to do it by name.
# Adds a "Names" worksheet (ID + Name only) positioned before "Products",
# pulls its data from the existing "Products" sheet, styles the Name
# column with wrap text + explicit black font, and tweaks a couple of
# selections left over on other sheets.

$wb = $excel.ActiveWorkbook

$productsSheet = $wb.Worksheets.Item("Products")
$cartSheet     = $wb.Worksheets.Item("Cart")
$inventorySheet = $wb.Worksheets.Item("Inventory")
$idSheet        = $wb.Worksheets.Item("ID")

# --- Create the new "Names" sheet right before "Products" ---------------
$namesSheet = $wb.Worksheets.Add($productsSheet)
$namesSheet.Name = "Names"

# --- Header + rows (ID, Name) pulled from the Products sheet ------------
$names = @(
    @("ID",  "Name"),
    @(4, "Sauce Labs Backpack"),
    @(0, "Sauce Labs Bike Light"),
    @(1, "Sauce Labs Bolt T-Shirt"),
    @(5, "Sauce Labs Fleece Jacket"),
    @(2, "Sauce Labs Onesie"),
    @(3, "Test.allTheThings() T-Shirt (Red)")
)

for ($i = 0; $i -lt $names.Length; $i++) {
    $r = $i + 1
    $namesSheet.Cells.Item($r, 1).Value = $names[$i][0]
    $namesSheet.Cells.Item($r, 2).Value = $names[$i][1]
}

# Column A keeps the plain wrap-text style already used for ID columns
# elsewhere in the workbook.
$namesSheet.Range("A1:A7").WrapText = $true

# Column B (Name) gets wrap text plus an explicit black font color,
# which mints the new font/cellXf pair.
$namesSheet.Range("B1:B7").WrapText = $true
$namesSheet.Range("B1:B7").Font.Color = 0

# Row heights to match the wrapped text layout.
$namesSheet.Rows(1).RowHeight = 17
$namesSheet.Rows(2).RowHeight = 34
$namesSheet.Rows(3).RowHeight = 34
$namesSheet.Rows(4).RowHeight = 34
$namesSheet.Rows(5).RowHeight = 51
$namesSheet.Rows(6).RowHeight = 34
$namesSheet.Rows(7).RowHeight = 51

# Selection on the new sheet.
$namesSheet.Range("A1:A7").Select()

# --- Leftover selection tweaks on other sheets ---------------------------
$inventorySheet.Range("C1:C7").Select()
$idSheet.Range("A1:A7").Select()

# --- Make "Names" the active tab (and thus the saved selection) ---------
$namesSheet.Activate()
